# "Started week 8 DQ1"
# Adds a new "week8" worksheet (after "week7"), seeded with the same
# task-tracker layout used by the other week sheets, with the first
# task (Discussion question 1 time estimate / actual / due date) filled in.

$wb = $excel.ActiveWorkbook
$ws7 = $wb.Worksheets.Item("week7")

# --- tidy up week7's view state: it stops being the tab-selected sheet ---
[void]$ws7.Select()
$excel.ActiveWindow.Zoom = 125
[void]$ws7.Range("B18").Select()

# --- create week8 right after week7 ---
$ws8 = $wb.Worksheets.Add($null, $ws7)
$ws8.Name = "week8"

# Pull over the column layout/formatting by copying whole rows from week7
# (this carries the number formats / styles used for every row) and then
# overwrite the copied values with week8's own content below.
$ws7.Range("A1:E1").Copy($ws8.Range("A1"))
for ($r = 2; $r -le 15; $r++) {
    $ws7.Range("A2:D2").Copy($ws8.Range("A$r"))
}
$ws7.Range("A14:D14").Copy($ws8.Range("A10"))
$ws7.Range("A14:D14").Copy($ws8.Range("A11"))
$ws7.Range("A4:D4").Copy($ws8.Range("A4"))
$ws7.Range("A18:D18").Copy($ws8.Range("A15"))
$ws7.Range("A19:D19").Copy($ws8.Range("A16"))

# --- row labels (task names), column A ---
$ws8.Range("A2").Value = "Discussion question 1"
$ws8.Range("A3").Value = "Discussion question 2"
$ws8.Range("A4").Value = "Read Lecture notes"
$ws8.Range("A5").Value = "DQ1 response 1"
$ws8.Range("A6").Value = "DQ1 response 2"
$ws8.Range("A7").Value = "DQ1 response 3"
$ws8.Range("A8").Value = "DQ1 response 4"
$ws8.Range("A9").Value = "DQ1 response 5"
$ws8.Range("A10").Value = "DQ2 response 1"
$ws8.Range("A11").Value = "DQ2 response 2"
$ws8.Range("A12").Value = "DQ2 response 3"
$ws8.Range("A13").Value = "DQ2 response 4"
$ws8.Range("A14").Value = "DQ2 response 5"
$ws8.Range("A15").Value = "Assignment"
$ws8.Range("A16").Value = "Total"

# --- anticipated length of time to complete, column B ---
$ws8.Range("B2").Value = 0.14583333333333334
$ws8.Range("B3").Value = 0.14583333333333334
$ws8.Range("B4").Value = 0.03125
$ws8.Range("B5").Value = 0.010416666666666666
$ws8.Range("B6").Value = 0.010416666666666666
$ws8.Range("B7").Value = 0.010416666666666666
$ws8.Range("B8").Value = 0.010416666666666666
$ws8.Range("B9").Value = 0.010416666666666666
$ws8.Range("B10").Value = 0.010416666666666666
$ws8.Range("B11").Value = 0.010416666666666666
$ws8.Range("B12").Value = 0.010416666666666666
$ws8.Range("B13").Value = 0.010416666666666666
$ws8.Range("B14").Value = 0.010416666666666666
$ws8.Range("B15").Value = 0.3333333333333333
$ws8.Range("B16").Formula = "=SUM(B2:B15)"

# --- actual time length to complete, column C ---
# Only DQ1 (row 2) has been started so far, 30 minutes in.
$ws8.Range("C2").Value = 0.020833333333333332
$ws8.Range("C3").ClearContents()
$ws8.Range("C4").ClearContents()
$ws8.Range("C5").ClearContents()
$ws8.Range("C6").ClearContents()
$ws8.Range("C7").ClearContents()
$ws8.Range("C8").ClearContents()
$ws8.Range("C9").ClearContents()
$ws8.Range("C10").ClearContents()
$ws8.Range("C11").ClearContents()
$ws8.Range("C12").ClearContents()
$ws8.Range("C13").ClearContents()
$ws8.Range("C14").ClearContents()
$ws8.Range("C15").ClearContents()
$ws8.Range("C16").Formula = "=SUM(C2:C15)"

# --- scheduled completion date, column D ---
$ws8.Range("D2").Value = 41824
$ws8.Range("D3").Value = 41825
$ws8.Range("D4").ClearContents()
$ws8.Range("D5").Value = 41825
$ws8.Range("D6").Value = 41826
$ws8.Range("D7").Value = 41827
$ws8.Range("D8").Formula = "=D7+1"
$ws8.Range("D9").Formula = "=D8+1"
$ws8.Range("D10").Value = 41825
$ws8.Range("D11").Value = 41826
$ws8.Range("D12").Value = 41827
$ws8.Range("D13").Formula = "=D12+1"
$ws8.Range("D14").Formula = "=D13+1"
$ws8.Range("D15").Value = 41829
$ws8.Range("D16").ClearContents()

# --- column widths (bestFit, same as the other week sheets) ---
$ws8.Columns.Item(1).ColumnWidth = 21.33203125
$ws8.Columns.Item(2).ColumnWidth = 32.83203125
$ws8.Columns.Item(3).ColumnWidth = 26.5
$ws8.Columns.Item(4).ColumnWidth = 23.5

# --- this is now the active, tab-selected sheet ---
[void]$ws8.Select()
$excel.ActiveWindow.Zoom = 125
[void]$ws8.Range("C3").Select()
